$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("nonsortable")

# Insert a new row before row 35, shifting rows 35-48 down to 36-49.
# Excel automatically extends the formatting of the surrounding rows into
# the freshly-inserted row, so no extra formatting work is required.
$ws.Rows.Item(35).Insert()

# Populate the new row's data cells.
$ws.Range("D35").Value = "Umpire review: Record keeping"
$ws.Range("E35").Value = "NK"

# Restore the autofilter / filter-database range to cover the extra row.
$ws.Range("A1:F43").AutoFilter()

$ws.Range("F35").Select()
